$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (GitHub Actions scheduled update)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.488.75"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.734.75"
$ws.Range("E3").Value = "  +4.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.37"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4782"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06216"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.735.99"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07129"
$ws.Range("E11").Value = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.69"
$ws.Range("E12").Value = "  +7.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6135"
$ws.Range("E13").Value = "  +7.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.522"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.70"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.504.55"
$ws.Range("E17").Value = "  +4.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006890"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.70"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.960.01"
$ws.Range("E21").Value = "  +4.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.556"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.876"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.319"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.60"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.795"
$ws.Range("E27").Value = "  +5.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.407"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.67"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.976"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.695"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07875"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04560"
$ws.Range("E33").Value = "  +5.32%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6326"
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9898"
$ws.Range("E36").Value = "  +5.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9304"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "111.01"
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.449"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.972"
$ws.Range("E40").Value = "  +8.29%  "
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01505"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("E43").Value = "  +13.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3891"
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.887"
$ws.Range("E45").Value = "  +12.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1190"
$ws.Range("E46").Value = "  +7.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05335"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.900"
$ws.Range("E48").Value = "  +4.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.73"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.250"
$ws.Range("E50").Value = "  +5.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3424"
$ws.Range("E51").Value = "  +3.63%  "
